# Update crypto price/volume data as scraped on Sat Aug 19 15:38:50 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    # Force the cell to remain a text value (matches the source sheet's
    # inline-string / shared-string cell typing) even when the text looks
    # like a plain number (e.g. "1.004"), then restore the default style
    # so no stray number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.144.87'
Set-TextValue $ws.Range("D3") '1.672.93'
Set-TextValue $ws.Range("E3") '  -0.13%  '
Set-TextValue $ws.Range("E4") '  -0.19%  '
Set-TextValue $ws.Range("D5") '216.82'
Set-TextValue $ws.Range("E5") '  -0.79%  '
Set-TextValue $ws.Range("D6") '0.5216'
Set-TextValue $ws.Range("E6") '  +1.94%  '
Set-TextValue $ws.Range("D7") '1.004'
Set-TextValue $ws.Range("E7") '  -0.18%  '
Set-TextValue $ws.Range("D8") '0.2698'
Set-TextValue $ws.Range("E8") '  +1.44%  '
Set-TextValue $ws.Range("D9") '0.06397'
Set-TextValue $ws.Range("E9") '  +0.99%  '
Set-TextValue $ws.Range("D10") '21.83'
Set-TextValue $ws.Range("E10") '  -0.30%  '
Set-TextValue $ws.Range("D11") '0.07429'
Set-TextValue $ws.Range("E11") '  +0.82%  '
Set-TextValue $ws.Range("D12") '1.697.21'
Set-TextValue $ws.Range("E12") '  +1.22%  '
Set-TextValue $ws.Range("D13") '4.520'
Set-TextValue $ws.Range("E13") '  -0.83%  '
Set-TextValue $ws.Range("D14") '0.5829'
Set-TextValue $ws.Range("E14") '  +1.12%  '
Set-TextValue $ws.Range("D15") '0.000008522'
Set-TextValue $ws.Range("E15") '  -0.10%  '
Set-TextValue $ws.Range("D16") '64.24'
Set-TextValue $ws.Range("E16") '  -1.34%  '
Set-TextValue $ws.Range("D17") '26.160.84'
Set-TextValue $ws.Range("E17") '  -0.60%  '
Set-TextValue $ws.Range("D18") '4.940'
Set-TextValue $ws.Range("E18") '  -1.37%  '
Set-TextValue $ws.Range("E19") '  -0.21%  '
Set-TextValue $ws.Range("D20") '10.79'
Set-TextValue $ws.Range("E20") '  -0.83%  '
Set-TextValue $ws.Range("D21") '189.65'
Set-TextValue $ws.Range("E21") '  +1.39%  '
Set-TextValue $ws.Range("E22") '  -0.46%  '
Set-TextValue $ws.Range("E23") '  -0.19%  '
Set-TextValue $ws.Range("D24") '144.83'
Set-TextValue $ws.Range("E24") '  +0.80%  '
Set-TextValue $ws.Range("B25") 'Cosmos'
Set-TextValue $ws.Range("C25") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D25") '7.619'
Set-TextValue $ws.Range("E25") '  +0.42%  '
Set-TextValue $ws.Range("B26") 'Stellar'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D26") '0.1239'
Set-TextValue $ws.Range("E26") '  +5.35%  '
Set-TextValue $ws.Range("D27") '15.73'
Set-TextValue $ws.Range("E27") '  +0.20%  '
Set-TextValue $ws.Range("D28") '0.06578'
Set-TextValue $ws.Range("E28") '  +12.99%  '
Set-TextValue $ws.Range("D29") '1.328'
Set-TextValue $ws.Range("E29") '  -0.10%  '
Set-TextValue $ws.Range("D30") '1.317'
Set-TextValue $ws.Range("E30") '  -0.54%  '
Set-TextValue $ws.Range("D31") '3.586'
Set-TextValue $ws.Range("E31") '  +2.08%  '
Set-TextValue $ws.Range("D32") '3.532'
Set-TextValue $ws.Range("E32") '  +0.75%  '
Set-TextValue $ws.Range("D33") '1.669'
Set-TextValue $ws.Range("E33") '  +0.30%  '
Set-TextValue $ws.Range("D34") '1.018'
Set-TextValue $ws.Range("E34") '  +1.45%  '
Set-TextValue $ws.Range("D35") '0.6167'
Set-TextValue $ws.Range("E35") '  +2.97%  '
Set-TextValue $ws.Range("D36") '2.367'
Set-TextValue $ws.Range("E36") '  +0.05%  '
Set-TextValue $ws.Range("D37") '2.699'
Set-TextValue $ws.Range("E37") '  +1.45%  '
Set-TextValue $ws.Range("D38") '6.273'
Set-TextValue $ws.Range("E38") '  +5.96%  '
Set-TextValue $ws.Range("D39") '1.095.99'
Set-TextValue $ws.Range("E39") '  -0.32%  '
Set-TextValue $ws.Range("E40") '  -0.62%  '
Set-TextValue $ws.Range("D41") '0.8723'
Set-TextValue $ws.Range("E41") '  +1.19%  '
Set-TextValue $ws.Range("E42") '  +0.63%  '
Set-TextValue $ws.Range("D43") '100.87'
Set-TextValue $ws.Range("E43") '  +1.41%  '
Set-TextValue $ws.Range("D44") '1.819.78'
Set-TextValue $ws.Range("E44") '  -0.32%  '
Set-TextValue $ws.Range("D45") '0.00000000109'
Set-TextValue $ws.Range("E45") '  -5.74%  '
Set-TextValue $ws.Range("E46") '  +0.29%  '
Set-TextValue $ws.Range("D47") '8.155'
Set-TextValue $ws.Range("E47") '  +1.10%  '
Set-TextValue $ws.Range("D48") '1.004'
Set-TextValue $ws.Range("E48") '  -0.33%  '
Set-TextValue $ws.Range("D49") '0.05241'
Set-TextValue $ws.Range("E49") '  +0.45%  '
Set-TextValue $ws.Range("D50") '0.4279'
Set-TextValue $ws.Range("E50") '  -0.82%  '
Set-TextValue $ws.Range("D51") '5.997'
Set-TextValue $ws.Range("E51") '  +2.75%  '
